$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-CellText "D2" "30.274.57"
Set-CellText "E2" "  +0.01%  "
Set-CellText "D3" "1.869.72"
Set-CellText "E3" "  -0.84%  "
Set-CellText "D5" "235.07"
Set-CellText "E5" "  -1.34%  "
Set-CellText "E6" "  +0.08%  "
Set-CellText "D7" "0.4660"
Set-CellText "E7" "  -0.29%  "
Set-CellText "D8" "0.2840"
Set-CellText "E8" "  +0.21%  "
Set-CellText "D9" "0.06557"
Set-CellText "E9" "  -0.40%  "
Set-CellText "D10" "21.21"
Set-CellText "E10" "  +6.72%  "
Set-CellText "D11" "0.07864"
Set-CellText "E11" "  +1.15%  "
Set-CellText "D12" "97.99"
Set-CellText "E12" "  -0.05%  "
Set-CellText "D13" "1.871.21"
Set-CellText "E13" "  -1.00%  "
Set-CellText "D14" "5.114"
Set-CellText "E14" "  -0.13%  "
Set-CellText "D15" "0.6750"
Set-CellText "E15" "  +0.89%  "
Set-CellText "D16" "281.30"
Set-CellText "E16" "  -0.84%  "
Set-CellText "D17" "30.266.93"
Set-CellText "E17" "  -0.02%  "
Set-CellText "D18" "1.000"
Set-CellText "E18" "  +0.02%  "
Set-CellText "D19" "5.522"
Set-CellText "E19" "  +2.89%  "
Set-CellText "D20" "12.68"
Set-CellText "E20" "  +0.53%  "
Set-CellText "D21" "2.116.85"
Set-CellText "E21" "  -0.79%  "
Set-CellText "D22" "0.000007288"
Set-CellText "E22" "  -0.22%  "
Set-CellText "E23" "  +0.04%  "
Set-CellText "D24" "6.172"
Set-CellText "E24" "  -0.05%  "
Set-CellText "D25" "9.230"
Set-CellText "E25" "  -1.29%  "
Set-CellText "D26" "164.60"
Set-CellText "E26" "  -0.23%  "
Set-CellText "D27" "19.25"
Set-CellText "E27" "  +0.39%  "
Set-CellText "D28" "1.929"
Set-CellText "E28" "  -3.23%  "
Set-CellText "D29" "1.373"
Set-CellText "E29" "  -0.15%  "
Set-CellText "D30" "0.09720"
Set-CellText "E30" "  -0.17%  "
Set-CellText "D31" "4.425"
Set-CellText "E31" "  -0.92%  "
Set-CellText "E32" "  -0.54%  "
Set-CellText "D33" "4.114"
Set-CellText "E33" "  -1.57%  "
Set-CellText "D34" "0.04689"
Set-CellText "E34" "  -0.13%  "
Set-CellText "D35" "1.122"
Set-CellText "E35" "  +2.41%  "
Set-CellText "D36" "0.7064"
Set-CellText "E36" "  -0.54%  "
Set-CellText "D37" "2.731"
Set-CellText "E37" "  +0.79%  "
Set-CellText "D38" "0.01857"
Set-CellText "E38" "  -0.81%  "
Set-CellText "D39" "6.278"
Set-CellText "E39" "  -5.54%  "
Set-CellText "D40" "2.540"
Set-CellText "E40" "  +0.57%  "
Set-CellText "D41" "73.59"
Set-CellText "E41" "  +1.57%  "
Set-CellText "D42" "1.950"
Set-CellText "E42" "  -1.18%  "
Set-CellText "D43" "0.8458"
Set-CellText "E43" "  -2.78%  "
Set-CellText "D44" "0.4175"
Set-CellText "E44" "  -0.68%  "
Set-CellText "B45" "PaxDollar"
Set-CellText "C45" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-CellText "D45" "1.001"
Set-CellText "E45" "  +0.11%  "
Set-CellText "B46" "Quant"
Set-CellText "C46" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-CellText "D46" "103.97"
Set-CellText "E46" "  -0.13%  "
Set-CellText "D47" "7.197"
Set-CellText "E47" "  -0.50%  "
Set-CellText "D48" "9.160"
Set-CellText "E48" "  -0.87%  "
Set-CellText "D49" "931.08"
Set-CellText "E49" "  -5.47%  "
Set-CellText "D50" "34.06"
Set-CellText "E50" "  -0.10%  "
Set-CellText "D51" "0.1125"
Set-CellText "E51" "  -3.21%  "
